$d = $word.ActiveDocument

# 1) Remove the hidden "_GoBack" bookmark that Word leaves behind from the
#    last editing position - it was cleaned up in this revision.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Insert a new run in red reading "//TEST" right after the
#    "...musi zmniejszyć czcionkę " text (small update of fontsize on the
#    card screen note).
$r = $d.Content
$found = $r.Find.Execute("czcionkę ", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
if ($found) {
    $insPos = $r.End
    $newRun = $d.Range($insPos, $insPos)
    $newRun.InsertAfter("//TEST")

    $newRunRange = $d.Range($insPos, $insPos + 6)
    $newRunRange.Font.Color = 255
}
